# Record the result of an "opt on eur/usd H1 2 years" run into Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 7: header / label row for the new block ---
$ws.Range("A7").Value = "opt on len/correl"
$ws.Range("B7").Value = "eur/usd"
$ws.Range("C7").Value = "H1"
$ws.Range("D7").Value = "10+"
$ws.Range("E7").Value = "90+"
$ws.Range("G7").Value = 30
$ws.Range("Q7").Value = "thougher conditions, better PF"

# --- Row 8 ---
$ws.Range("B8").Value = "eur/usd"
$ws.Range("C8").Value = "H1"
$ws.Range("D8").Value = 8
$ws.Range("E8").Value = 92
$ws.Range("G8").Value = 30
$ws.Range("H8").Value = 25
$ws.Range("J8").Value = 10000
$ws.Range("L2").Copy()
$ws.Range("L8").PasteSpecial(-4122)  # xlPasteFormats, match L2's date format/style
$ws.Range("L8").Value = 41275        # 1/1/2013 (days since 1899-12-30)
$ws.Range("N8").Value = 100
$ws.Range("O8").Value = 2100
$ws.Range("P8").Value = 1.7
$ws.Range("Q8").Value = "to 1/4/2015"

# --- Row 9 ---
$ws.Range("B9").Value = "eur/usd"
$ws.Range("C9").Value = "H1"
$ws.Range("D9").Value = 18
$ws.Range("E9").Value = 86
$ws.Range("N9").Value = 62
$ws.Range("O9").Value = 2000
$ws.Range("P9").Value = 2.4

# --- Row 10 ---
$ws.Range("D10").Value = 14
$ws.Range("E10").Value = 88
$ws.Range("N10").Value = 82
$ws.Range("O10").Value = 1900
$ws.Range("P10").Value = 1.7

# --- Row 11 ---
$ws.Range("D11").Value = 20
$ws.Range("E11").Value = 86
$ws.Range("N11").Value = 50
$ws.Range("O11").Value = 1700
$ws.Range("P11").Value = 2

# --- Row 12 ---
$ws.Range("D12").Value = 18
$ws.Range("E12").Value = 94
$ws.Range("N12").Value = 10
$ws.Range("O12").Value = 1000
$ws.Range("P12").Value = 7.5

# Apply the same fill styles used by the existing D/E columns (style index 2
# in the raw XML, i.e. the light highlight fill) to the new D7:E12 cells.
$ws.Range("D2:E2").Copy()
$ws.Range("D7:E12").PasteSpecial(-4122)  # xlPasteFormats

# Widen column Q to fit the new, longer note text (closest achievable width
# to the recorded 69.7109375 given this host's column-width quantization).
$ws.Columns.Item(17).ColumnWidth = 68.8

# Update the active selection to match where the user left off after typing.
$ws.Range("P13").Select()
